$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from 45203 to 45204 for rows 2..349
$ws.Range("C2:C349").Value = 45204

# Add new row 350 with data
$ws.Range("A350").Value = "A 47604-2023"
$ws.Range("B350").Value = 45203
$ws.Range("C350").Value = 45204
$ws.Range("D350").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E350").Value = "NORSJÖ"
$ws.Range("F350").Value = "Holmen skog AB"
$ws.Range("G350").Value = 1.7
$ws.Range("H350").Value = 0
$ws.Range("I350").Value = 0
$ws.Range("J350").Value = 0
$ws.Range("K350").Value = 0
$ws.Range("L350").Value = 0
$ws.Range("M350").Value = 0
$ws.Range("N350").Value = 0
$ws.Range("O350").Value = 0
$ws.Range("P350").Value = 0
$ws.Range("Q350").Value = 0

# Copy style from B349/C349 to B350/C350 (date format)
$ws.Range("B349").Copy()
$ws.Range("B350").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("C349").Copy()
$ws.Range("C350").PasteSpecial(-4122)

$ws.Range("R349").Copy()
$ws.Range("R350").PasteSpecial(-4122)

# Explicitly set row height for row 349 (matches target where it becomes explicit)
$ws.Rows.Item(349).RowHeight = 15
